$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.00"
$ws.Range("E2").Value = "'-0.35%"
$ws.Range("G2").Value = "'4"
$ws.Range("D3").Value = "'26.55"
$ws.Range("E3").Value = "'4.65%"
$ws.Range("G3").Value = "'4"
$ws.Range("D4").Value = "'5.123"
$ws.Range("E4").Value = "'-0.41%"
$ws.Range("G4").Value = "'4"
$ws.Range("D5").Value = "'0.05589"
$ws.Range("E5").Value = "'0.25%"
$ws.Range("G5").Value = "'4"
$ws.Range("D6").Value = "'6.472"
$ws.Range("E6").Value = "'-0.57%"
$ws.Range("G6").Value = "'4"
$ws.Range("D7").Value = "'0.8169"
$ws.Range("E7").Value = "'-0.21%"
$ws.Range("G7").Value = "'4"
$ws.Range("D8").Value = "'0.8355"
$ws.Range("E8").Value = "'-1.63%"
$ws.Range("G8").Value = "'4"
$ws.Range("D9").Value = "'0.1333"
$ws.Range("E9").Value = "'-0.69%"
$ws.Range("G9").Value = "'4"
$ws.Range("D10").Value = "'0.06994"
$ws.Range("E10").Value = "'0.49%"
$ws.Range("G10").Value = "'4"
$ws.Range("E11").Value = "'0.24%"
$ws.Range("G11").Value = "'4"
$ws.Range("D12").Value = "'0.09395"
$ws.Range("E12").Value = "'0.15%"
$ws.Range("G12").Value = "'4"
$ws.Range("D13").Value = "'0.001528"
$ws.Range("E13").Value = "'0.79%"
$ws.Range("G13").Value = "'4"
$ws.Range("E14").Value = "'0.60%"
$ws.Range("G14").Value = "'4"
$ws.Range("D15").Value = "'0.006170"
$ws.Range("E15").Value = "'1.14%"
$ws.Range("G15").Value = "'4"
$ws.Range("D16").Value = "'3.647"
$ws.Range("E16").Value = "'4.17%"
$ws.Range("G16").Value = "'4"
$ws.Range("D17").Value = "'3.038"
$ws.Range("E17").Value = "'0.63%"
$ws.Range("G17").Value = "'4"
$ws.Range("D18").Value = "'2.182"
$ws.Range("E18").Value = "'5.75%"
$ws.Range("G18").Value = "'4"
$ws.Range("G19").Value = "'4"
$ws.Range("D20").Value = "'0.03126"
$ws.Range("E20").Value = "'-1.23%"
$ws.Range("G20").Value = "'4"
$ws.Range("E21").Value = "'-2.26%"
$ws.Range("G21").Value = "'4"
$ws.Range("D22").Value = "'3.758"
$ws.Range("E22").Value = "'-0.02%"
$ws.Range("G22").Value = "'4"
$ws.Range("D23").Value = "'0.04576"
$ws.Range("E23").Value = "'-3.29%"
$ws.Range("G23").Value = "'4"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("G24").Value = "'4"
$ws.Range("D25").Value = "'0.001244"
$ws.Range("E25").Value = "'-0.33%"
$ws.Range("G25").Value = "'4"
$ws.Range("D26").Value = "'0.004506"
$ws.Range("E26").Value = "'-2.87%"
$ws.Range("G26").Value = "'4"
$ws.Range("D27").Value = "'0.00009601"
$ws.Range("E27").Value = "'-1.05%"
$ws.Range("G27").Value = "'4"
$ws.Range("E28").Value = "'0.51%"
$ws.Range("G28").Value = "'4"
$ws.Range("G29").Value = "'4"
$ws.Range("G30").Value = "'4"
$ws.Range("G31").Value = "'4"
$ws.Range("G32").Value = "'4"
$ws.Range("G33").Value = "'4"
$ws.Range("G34").Value = "'4"
$ws.Range("G35").Value = "'4"
$ws.Range("G36").Value = "'4"
$ws.Range("G37").Value = "'4"
$ws.Range("G38").Value = "'4"
$ws.Range("G39").Value = "'4"
$ws.Range("G40").Value = "'4"
$ws.Range("D41").Value = "'0.1379"
$ws.Range("E41").Value = "'31.28%"
$ws.Range("G41").Value = "'4"
$ws.Range("D42").Value = "'0.006161"
$ws.Range("E42").Value = "'-0.50%"
$ws.Range("G42").Value = "'4"
$ws.Range("D43").Value = "'0.002590"
$ws.Range("E43").Value = "'-1.32%"
$ws.Range("G43").Value = "'4"
$ws.Range("D44").Value = "'0.008854"
$ws.Range("E44").Value = "'6.67%"
$ws.Range("G44").Value = "'4"
$ws.Range("D45").Value = "'0.00005333"
$ws.Range("E45").Value = "'0.62%"
$ws.Range("G45").Value = "'4"
$ws.Range("G46").Value = "'4"
$ws.Range("D47").Value = "'0.1440"
$ws.Range("E47").Value = "'-23.87%"
$ws.Range("G47").Value = "'4"
$ws.Range("D48").Value = "'0.002337"
$ws.Range("E48").Value = "'10.17%"
$ws.Range("G48").Value = "'4"
$ws.Range("G49").Value = "'4"
$ws.Range("G50").Value = "'4"
$ws.Range("G51").Value = "'4"
